$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was captured and inserted as row 14 (pushing the
# previously-existing rows 14-40 down to 15-41). Insert a fresh row at 14
# so everything below shifts down by one, matching Excel's native
# "insert row" behaviour (formats/row-shift handled by the host).
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new observation.
$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "Macroferia Regional de Talca"
$ws.Range("C14").Value = "Maule"
$ws.Range("D14").Value = 44665
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100104
$ws.Range("H14").Value = "Frutos de pepita"
$ws.Range("I14").Value = 100104003
$ws.Range("J14").Value = "Membrillo"
$ws.Range("K14").Value = "Champion"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 250
$ws.Range("N14").Value = 10000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 10000
$ws.Range("Q14").Value = '$/caja 18 kilos granel'
$ws.Range("R14").Value = "Región de O'Higgins"
$ws.Range("S14").Value = 556
$ws.Range("T14").Value = 18
